# Update the "Förändrad" (Changed) date column (C) for rows 2-14.
# The stored date serial number moves from 46075 (2026-02-22) to
# 46076 (2026-02-23) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
